$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update the "Marking" row (right-answer marking value) and the "Total" row
# (total correct marks + the Corr/total marks summary string).
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 105
$ws.Range("E12").Value = "105/140"
